# Atualização automática OEE, Canudos, Producao e Rejeito (06/02/2026 14:10:13,65)
# Adds three new daily-production log rows (108-110) to the "Canudos" log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 108: OS 8502, turno B, 6179 produzidos, 25 perdas, cor 475 ---
$ws.Range("A108").Value = 46056
$ws.Range("B108").Value = "B"
$ws.Range("D108").Value = 8502
$ws.Range("E108").Value = 6179
$ws.Range("F108").Value = 25
$ws.Range("G108").Value = "cor 475"

# --- Row 109: OS 8524, turno A, problema na queda de canudos ---
$ws.Range("A109").Value = 46057
$ws.Range("B109").Value = "A"
$ws.Range("D109").Value = 8524
$ws.Range("F109").Value = " "
$ws.Range("G109").Value = "MAQUINA APRESENTOU PROBLEMA NA QUEDA DE CANUDOS  DA CALHA PARA A CORRENTE, TENTIVA DE MANUTENÇÃO POR OPERADOR RESPONSAVEL"
$ws.Range("H109").Value = " "

# --- Row 110: OS 8524, turno A, automação indisponível ---
$ws.Range("A110").Value = 46058
$ws.Range("B110").Value = "A"
$ws.Range("D110").Value = 8524
$ws.Range("F110").Value = " "
$ws.Range("G110").Value = "PROBLEMA ENCONTRADO, PORÉM, APENAS AUTOMAÇÃO PODE RESOLVER, OPERDOR LIGOU PARA O RESPONSAVEL DA MAQUINA NA AUTOMAÇÃO MAS O MESMO ESTAVA AFASTADO NO DIA POR ATESDADO"
$ws.Range("H110").Value = " "
$ws.Range("I110").Value = " "

# Scroll/selection update to reflect the newly added rows being in view.
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J120").Select()

$wb.Save()
